$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "Message for New transmittal" value for the Multi user record (row 3),
# matching the same text already present in row 2's L column.
$ws.Range("L3").Value = "Message for New transmittal"

# Reflect the new active selection on the sheet (matches the saved view state).
$ws.Range("L3").Select()
